$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Project description paragraph: replace placeholder text with the
#    real description, split across two runs (same formatting) exactly
#    like the authored edit (Word naturally breaks runs while typing).
# ---------------------------------------------------------------------
$rng1 = $d.Content
$old1 = "<Scrivete qua una breve descrizione del vostro progetto, max 5 righe>"
$firstPart = "Lo scopo del progetto è arricchire la vita degli studenti del politecnico di Milano s"
$secondPart = "otto il profilo sociale e sportivo"
$found1 = $rng1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $firstPart, 2)
if ($found1) {
    $insAfter = $d.Range($rng1.End, $rng1.End)
    $insAfter.InsertAfter($secondPart)
    $splitRng = $d.Range($rng1.End, $rng1.End + $secondPart.Length)
    # Force a real run boundary (Word keeps the two runs distinct even
    # though the resulting formatting is identical) by toggling a
    # character property back to its original value.
    $splitRng.Bold = 1
    $splitRng.Bold = 0
}

# ---------------------------------------------------------------------
# 2) License paragraph: move the "diritto di riproduzione ... adattamento, "
#    sentence fragment from the end of the first run to the start of the
#    second run (which carries a <w:lastRenderedPageBreak/>). Trim it off
#    the first run, then insert it right before the break run's text so
#    the page-break marker stays attached to "elaborazione e riduzione...".
# ---------------------------------------------------------------------
$movedText = "diritto di riproduzione in qualunque modo o forma; diritto di trascrizione, montaggio, adattamento, "
$rngTrim = $d.Content
$oldTrim = "diritto di pubblicazione; " + $movedText
$newTrim = "diritto di pubblicazione; "
$foundTrim = $rngTrim.Find.Execute($oldTrim, $false, $false, $false, $false, $false, $true, 1, $false, $newTrim, 2)

if ($foundTrim) {
    $rngAfterBreak = $d.Content
    $foundBreak = $rngAfterBreak.Find.Execute("elaborazione e riduzione; diritto di comunicazione", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($foundBreak) {
        $insPoint = $d.Range($rngAfterBreak.Start, $rngAfterBreak.Start)
        $insPoint.InsertBefore($movedText)
    }
}

# ---------------------------------------------------------------------
# 3) Drop the stray <w:lastRenderedPageBreak/> in front of the second
#    "CONSENSO AL TRATTAMENTO DEI DATI PERSONALI" heading (first
#    occurrence, inside "DICHIARAZIONE DI CONSENSO ...", must stay
#    untouched). Re-finding/replacing the run's own text with itself
#    regenerates the run without the stale break marker.
# ---------------------------------------------------------------------
$rngSkip = $d.Content
$foundSkip = $rngSkip.Find.Execute("AI SENSI DELL'ART. 13 DEL REGOLAMENTO", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundSkip) {
    $rngHeading = $d.Range($rngSkip.End, $d.Content.End)
    $headingText = "CONSENSO AL TRATTAMENTO DEI DATI PERSONALI "
    $foundHeading = $rngHeading.Find.Execute($headingText, $false, $false, $false, $false, $false, $true, 1, $false, $headingText, 2)
}

# ---------------------------------------------------------------------
# 4) "nato/a a ......." : collapse the "nato/a " + "a" (wrapped in
#    proofErr spell-check markers) + " ...... il......" runs into one
#    plain run with identical text.
# ---------------------------------------------------------------------
$rng4 = $d.Content
$text4 = "nato/a a …………………… il………………………… residente a …………………………………………… "
$found4 = $rng4.Find.Execute($text4, $false, $false, $false, $false, $false, $true, 1, $false, $text4, 2)

# ---------------------------------------------------------------------
# 5) "via..................." : collapse the address-line runs (split up
#    by grammar-check proofErr markers) into a single run with identical
#    text.
# ---------------------------------------------------------------------
$rng5 = $d.Content
$text5 = "via………………………………………………………………….. n……… ; "
$found5 = $rng5.Find.Execute($text5, $false, $false, $false, $false, $false, $true, 1, $false, $text5, 2)

Write-Output "edit1=$found1 trim=$foundTrim break=$foundBreak heading=$foundHeading nato=$found4 via=$found5"
